$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("my_products")

# Update the product list value in B2 (simplify "Credit Cards: Silver, Gold, Platinum & Millennium" to "Credit Cards")
$ws.Range("B2").Value = "Credit Cards|Deposit Account|Mortgage Home Loans|Personal Loans|Savings Account|Wealth Management"

# Move the selection from B10 to B9
$ws.Range("B9").Select()
